# Update "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet, and add a newly-scraped con ("南昌·LY-COSPLAY大会X运动
# 番PRO2.0（非ONLY）") as a new row 36 on the "全部类型" sheet (pushing the
# rows that used to be 36-38 down to 37-39).
#
# NOTE: reading a Range's `.Value` getter is broken in this host (it hands
# back the property's reflection signature instead of the cell content), so
# whenever we need to read a cell back we use `.Value2` instead, which works
# correctly for both numbers and (unicode) strings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")

$exhibitUpdates = @{
    "F2"  = 209
    "F3"  = 5551
    "F6"  = 33
    "F7"  = 661
    "F9"  = 8
    "F12" = 1549
    "F13" = 5104
    "F15" = 243
    "F16" = 215
    "F17" = 29
    "F18" = 12
    "F20" = 4402
    "F21" = 211
    "F22" = 1160
    "F27" = 175
    "F28" = 63
    "F29" = 148
    "F31" = 344
    "F33" = 42
}

foreach ($addr in $exhibitUpdates.Keys) {
    $wsExhibit.Range($addr).Value2 = $exhibitUpdates[$addr]
}

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    "F2"  = 209
    "F4"  = 5551
    "F7"  = 33
    "F8"  = 661
    "F10" = 8
    "F13" = 1549
    "F14" = 5104
    "F16" = 243
    "F17" = 215
    "F18" = 29
    "F19" = 12
    "F21" = 4402
    "F22" = 211
    "F23" = 1160
    "F28" = 175
    "F29" = 63
    "F30" = 148
    "F32" = 344
    "F34" = 42
}

foreach ($addr in $allUpdates.Keys) {
    $wsAll.Range($addr).Value2 = $allUpdates[$addr]
}

# Make room for the new con by sliding the data (columns B..I) of rows
# 36-38 down into rows 37-39; column A (the running index) is left exactly
# as it is for every row that already existed, and only gets a brand-new
# value for the row that is newly appearing at the bottom (row 39).
foreach ($pair in @(@(38, 39), @(37, 38), @(36, 37))) {
    $src = $pair[0]
    $dst = $pair[1]
    foreach ($col in @("C", "D", "E", "F", "G", "H", "I")) {
        $wsAll.Range("$col$dst").Value2 = $wsAll.Range("$col$src").Value2
    }
    # Column B holds plain-text dates ("YYYY-MM-DD"); force text formatting
    # on the destination so Excel doesn't reinterpret it as a date serial.
    $wsAll.Range("B$dst").NumberFormat = "@"
    $wsAll.Range("B$dst").Value2 = $wsAll.Range("B$src").Value2
}

# Row 39 is brand new, so its "A" cell needs to be created (matching the
# bordered/bold look of the rest of the index column) and given the next
# sequential value.
$wsAll.Range("A38").EntireRow.Copy() | Out-Null
$wsAll.Range("A39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$wsAll.Range("A39").Value2 = 38

# Finally, fill row 36 with the newly-scraped con.
$wsAll.Range("B36").NumberFormat = "@"
$wsAll.Range("B36").Value2 = "2024-06-10"
$wsAll.Range("C36").Value2 = "南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）"
$wsAll.Range("D36").Value2 = "青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK"
$wsAll.Range("E36").Value2 = "2024.06.10 10:00-06.10 17:00"
$wsAll.Range("F36").Value2 = 0
$wsAll.Range("G36").Value2 = 30
$wsAll.Range("H36").Value2 = "https://show.bilibili.com/platform/detail.html?id=84575"
$wsAll.Range("I36").Value2 = "//i2.hdslb.com/bfs/openplatform/202404/ScwkijwU1713428452963.jpeg"

$wb.Save()
